$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$codes = @(
    "2Let",
    "CN",
    "J",
    "KR",
    "IND",
    "MY",
    "SG",
    "INDO",
    "PH",
    "TH",
    "VN",
    "CA",
    "US",
    "MX",
    "BR",
    "AR",
    "B",
    "F",
    "D",
    "I",
    "NL",
    "CH",
    "E",
    "GB",
    "A",
    "GR",
    "S",
    "RUS",
    "AU",
    "NZ",
    "ZA"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $codes[$i]
}

$ws.Range("C1:C31").Select() | Out-Null
